# New crime data collected - weekly CompStat update
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Header text updates (Volume/Number banner + "Report covering the week" banner)
# ---------------------------------------------------------------------------
$ws.Range("A8").Value = "Volume 31   Number  12"
$ws.Range("C9").Value = "Report Covering the Week  3/18/2024  Through  3/24/2024"

# ---------------------------------------------------------------------------
# Helper: convert a numeric-looking cell to a genuine text cell (shared
# string) while keeping the cell's existing "text" look (font/format copied
# from a donor cell that already uses the text style used throughout this
# sheet for "0" / "***.*" placeholders).
# ---------------------------------------------------------------------------
function Set-TextCell($addr, $donor, $text) {
    $ws.Range($addr).Value = "'" + $text
    $ws.Range($donor).Copy() | Out-Null
    $ws.Range($addr).PasteSpecial(-4122) | Out-Null
}

# Helper: convert a text cell back into a genuine numeric cell while
# re-using the numeric style already used by a sibling numeric cell.
function Set-NumberCell($addr, $donor, $num) {
    $ws.Range($donor).Copy() | Out-Null
    $ws.Range($addr).PasteSpecial(-4122) | Out-Null
    $ws.Range($addr).Value = $num
}

# ---------------------------------------------------------------------------
# Row 16 - Robbery
# ---------------------------------------------------------------------------
$ws.Range("F16").Value = 3
$ws.Range("G16").Value = 6
$ws.Range("H16").Value = -50
$ws.Range("J16").Value = 12
$ws.Range("K16").Value = -33.333333333333
$ws.Range("N16").Value = -86.206896551724

# ---------------------------------------------------------------------------
# Row 17 - Fel. Assault
# ---------------------------------------------------------------------------
$ws.Range("C17").Value = 5
Set-TextCell "D17" "C14" "0"
Set-TextCell "E17" "E14" "***.*"
$ws.Range("F17").Value = 16
$ws.Range("H17").Value = 220
$ws.Range("I17").Value = 36
$ws.Range("K17").Value = 89.473684210526
$ws.Range("L17").Value = 89.473684210526
$ws.Range("M17").Value = 100
$ws.Range("N17").Value = -23.404255319148

# ---------------------------------------------------------------------------
# Row 18 - Burglary
# ---------------------------------------------------------------------------
$ws.Range("C18").Value = 2
Set-TextCell "D18" "C14" "0"
Set-TextCell "E18" "E14" "***.*"
$ws.Range("F18").Value = 6
$ws.Range("G18").Value = 1
$ws.Range("H18").Value = 500
$ws.Range("I18").Value = 11
$ws.Range("K18").Value = 57.142857142857
$ws.Range("L18").Value = 10
$ws.Range("M18").Value = -15.384615384615
$ws.Range("N18").Value = -91.2

# ---------------------------------------------------------------------------
# Row 19 - Gr. Larceny
# ---------------------------------------------------------------------------
$ws.Range("D19").Value = 7
$ws.Range("E19").Value = -71.428571428571
$ws.Range("F19").Value = 10
$ws.Range("G19").Value = 13
$ws.Range("H19").Value = -23.076923076923
$ws.Range("I19").Value = 32
$ws.Range("J19").Value = 40
$ws.Range("K19").Value = -20
$ws.Range("L19").Value = 0
$ws.Range("M19").Value = 28
$ws.Range("N19").Value = 6.666666666666

# ---------------------------------------------------------------------------
# Row 20 - G.L.A.
# ---------------------------------------------------------------------------
$ws.Range("C20").Value = 2
$ws.Range("F20").Value = 4
$ws.Range("H20").Value = 33.333333333333
$ws.Range("I20").Value = 17
$ws.Range("K20").Value = 142.857142857143
$ws.Range("L20").Value = 70
$ws.Range("M20").Value = 240
$ws.Range("N20").Value = -70.689655172413

# ---------------------------------------------------------------------------
# Row 21 - TOTAL
# ---------------------------------------------------------------------------
$ws.Range("C21").Value = 11
$ws.Range("D21").Value = 9
$ws.Range("E21").Value = 22.222222222222
$ws.Range("F21").Value = 39
$ws.Range("G21").Value = 28
$ws.Range("H21").Value = 39.285714285714
$ws.Range("I21").Value = 104
$ws.Range("J21").Value = 86
$ws.Range("K21").Value = 20.930232558139
$ws.Range("L21").Value = 28.395061728395
$ws.Range("M21").Value = 36.842105263157
$ws.Range("N21").Value = -68

# ---------------------------------------------------------------------------
# Row 22 - Transit
# ---------------------------------------------------------------------------
Set-NumberCell "C22" "I23" 2
Set-TextCell "D22" "C14" "0"
Set-TextCell "E22" "E14" "***.*"
$ws.Range("F22").Value = 4
$ws.Range("H22").Value = 300
$ws.Range("I22").Value = 5
$ws.Range("K22").Value = 400
$ws.Range("L22").Value = 400

# ---------------------------------------------------------------------------
# Row 23 - Housing
# ---------------------------------------------------------------------------
Set-NumberCell "C23" "I23" 1
Set-NumberCell "D23" "I23" 1
Set-NumberCell "E23" "K23" 0
Set-NumberCell "F23" "I23" 1
Set-NumberCell "G23" "I23" 1
Set-NumberCell "H23" "K23" 0
$ws.Range("I23").Value = 5
$ws.Range("J23").Value = 6
$ws.Range("K23").Value = -16.666666666666
$ws.Range("L23").Value = 25
$ws.Range("M23").Value = -44.444444444444

# ---------------------------------------------------------------------------
# Row 24 - Petit Larceny
# ---------------------------------------------------------------------------
$ws.Range("C24").Value = 12
$ws.Range("D24").Value = 12
$ws.Range("E24").Value = 0
$ws.Range("F24").Value = 50
$ws.Range("H24").Value = 35.135135135135
$ws.Range("I24").Value = 102
$ws.Range("J24").Value = 120
$ws.Range("K24").Value = -15
$ws.Range("L24").Value = 13.333333333333
$ws.Range("M24").Value = 61.904761904761

# ---------------------------------------------------------------------------
# Row 25 - Retail Theft
# ---------------------------------------------------------------------------
$ws.Range("C25").Value = 3
$ws.Range("D25").Value = 11
$ws.Range("E25").Value = -72.727272727272
$ws.Range("G25").Value = 31
$ws.Range("H25").Value = -45.161290322580
$ws.Range("I25").Value = 39
$ws.Range("J25").Value = 78
$ws.Range("K25").Value = -50
$ws.Range("L25").Value = 5.405405405405

# ---------------------------------------------------------------------------
# Row 26 - Misd. Assault
# ---------------------------------------------------------------------------
$ws.Range("C26").Value = 8
$ws.Range("D26").Value = 8
$ws.Range("E26").Value = 0
$ws.Range("F26").Value = 20
$ws.Range("G26").Value = 14
$ws.Range("H26").Value = 42.857142857142
$ws.Range("I26").Value = 50
$ws.Range("J26").Value = 43
$ws.Range("K26").Value = 16.279069767441
$ws.Range("L26").Value = 38.888888888888
$ws.Range("M26").Value = -20.634920634920

# ---------------------------------------------------------------------------
# Row 28 - Other Sex Crimes
# ---------------------------------------------------------------------------
Set-NumberCell "C28" "I23" 1
Set-TextCell "D28" "C14" "0"
Set-TextCell "E28" "E14" "***.*"
Set-NumberCell "F28" "I23" 1
$ws.Range("H28").Value = -50
$ws.Range("I28").Value = 6
$ws.Range("K28").Value = 20
$ws.Range("L28").Value = 0

# ---------------------------------------------------------------------------
# Row 29 - Shooting Vic.
# ---------------------------------------------------------------------------
Set-TextCell "C29" "C14" "0"

# ---------------------------------------------------------------------------
# Row 30 - Shooting Inc.
# ---------------------------------------------------------------------------
Set-TextCell "C30" "C14" "0"

# ---------------------------------------------------------------------------
# Row 31 - Hate Crimes
# ---------------------------------------------------------------------------
Set-NumberCell "C31" "I23" 1
Set-NumberCell "F31" "I23" 1
Set-NumberCell "I31" "I23" 1

# ---------------------------------------------------------------------------
# Row 33 - Traffic Fatalities
# ---------------------------------------------------------------------------
Set-NumberCell "C33" "I23" 1
Set-NumberCell "F33" "I23" 1
Set-NumberCell "I33" "I23" 1

Write-Host "edits applied"
